# "Generate Report for Handback"
# Refresh the timestamps recorded in the handback-status report:
#  - Overview sheet: "Latest HO Xliff Generate Date" (also shared by the
#    de-de sheet's "Correspond Handoff Datetime" for the same file, since
#    both cells hold the identical original timestamp string).
#  - zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
#  - de-de sheet: "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 and de-de!H2 both currently hold "2016-08-27 15:02:48";
# update both to the newly generated timestamp.
$wsOverview.Range("G2").Value = "2016-08-27 15:03:32"
$wsDeDe.Range("H2").Value = "2016-08-27 15:03:32"

# zh-cn!H2 "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-27 15:03:26"

# zh-cn!K2 "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-08-27 15:03:49"

# de-de!K2 "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-08-27 15:03:56"
